# 直播源汇总文档/其他文档/webview频道统计.xlsx
# "广西" (Guangxi) sheet: add three new Liuzhou (柳州) channel rows
# (new live.lzgd.com.cn stream URLs for the existing 柳州新闻综合 channel,
# plus a brand-new 柳州科教 channel with both a live.lzgd.com.cn and a
# web.guangdianyun.tv source), inserted right above the pre-existing
# 柳州新闻综合 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("广西")

# Make room for three new rows above the existing row 11 (柳州新闻综合),
# pushing the old row 11 down to row 14 along with everything below it.
$ws.Range("A11:A13").Insert()

# Row 11: new 柳州新闻综合 source (live.lzgd.com.cn)
$ws.Range("A11").Value = "柳州新闻综合,webview://https://live.lzgd.com.cn/tv/1902?uin=3370"
# Row 12: the original 柳州新闻综合 source (web.guangdianyun.tv), unchanged text
$ws.Range("A12").Value = "柳州新闻综合,webview://https://web.guangdianyun.tv/tv/?id=1902&uin=3370"
# Row 13: new 柳州科教 source (live.lzgd.com.cn)
$ws.Range("A13").Value = "柳州科教,webview://https://live.lzgd.com.cn/tv/1808?uin=3370"
# Row 14: new 柳州科教 source (web.guangdianyun.tv)
$ws.Range("A14").Value = "柳州科教,webview://https://web.guangdianyun.tv/tv/?id=1808&uin=3370"
